$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("pages_with_public_use")
$ws2 = $wb.Worksheets.Item("pages_without_public_use")

# Sheet1: update row4 path and delete row5 (old "espanol/esp-press-releae-no-date")
$ws1.Range("A4").Value = "news-events/press-releases/2018/leukemia-cll-ibrutinib-trial"
$ws1.Rows.Item(5).Delete()

# Sheet2: update row5 path/type
$ws2.Range("A5").Value = "about-cancer/coping/feelings/relaxation/loukissas-jennifer"
$ws2.Range("B5").Value = "Bio"

# Column D on sheet1 needs to widen to fit the new (longer) path text;
# 47.666666666666664 round-trips to the OOXML width value 48.5
$ws1.Columns.Item(4).ColumnWidth = 47.666666666666664

# Restore view/selection state on both sheets
$ws2.Activate()
[void]$ws2.Range("A10").Select()
$ws1.Activate()
[void]$ws1.Range("A13").Select()

